$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-10-08 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-09 Wednesday", 2) | Out-Null

# Update the math-facts table cells in row-major order
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "19+77="
$t.Cell(1, 2).Range.Text = "81-67="
$t.Cell(1, 3).Range.Text = "2+2="
$t.Cell(1, 4).Range.Text = "27+32="
$t.Cell(1, 5).Range.Text = "44+2="
$t.Cell(2, 1).Range.Text = "89-79="
$t.Cell(2, 2).Range.Text = "62+36="
$t.Cell(2, 3).Range.Text = "69-17="
$t.Cell(2, 4).Range.Text = "65-65="
$t.Cell(2, 5).Range.Text = "96-89="
$t.Cell(3, 1).Range.Text = "98-90="
$t.Cell(3, 2).Range.Text = "1+87="
$t.Cell(3, 3).Range.Text = "74-57="
$t.Cell(3, 4).Range.Text = "51-8="
$t.Cell(3, 5).Range.Text = "80-59="
$t.Cell(4, 1).Range.Text = "82-45="
$t.Cell(4, 2).Range.Text = "48+16="
$t.Cell(4, 3).Range.Text = "37+22="
$t.Cell(4, 4).Range.Text = "32+44="
$t.Cell(4, 5).Range.Text = "27+38="
$t.Cell(5, 1).Range.Text = "24+25="
$t.Cell(5, 2).Range.Text = "98-78="
$t.Cell(5, 3).Range.Text = "56-5="
$t.Cell(5, 4).Range.Text = "26+51="
$t.Cell(5, 5).Range.Text = "48+20="
$t.Cell(6, 1).Range.Text = "86-7="
$t.Cell(6, 2).Range.Text = "91-67="
$t.Cell(6, 3).Range.Text = "44+54="
$t.Cell(6, 4).Range.Text = "64-27="
$t.Cell(6, 5).Range.Text = "46+18="
$t.Cell(7, 1).Range.Text = "51+5="
$t.Cell(7, 2).Range.Text = "40-26="
$t.Cell(7, 3).Range.Text = "26+50="
$t.Cell(7, 4).Range.Text = "0+60="
$t.Cell(7, 5).Range.Text = "44-28="
$t.Cell(8, 1).Range.Text = "55-27="
$t.Cell(8, 2).Range.Text = "39-19="
$t.Cell(8, 3).Range.Text = "76+14="
$t.Cell(8, 4).Range.Text = "78-53="
$t.Cell(8, 5).Range.Text = "56-17="
$t.Cell(9, 1).Range.Text = "66-56="
$t.Cell(9, 2).Range.Text = "62-18="
$t.Cell(9, 3).Range.Text = "31+17="
$t.Cell(9, 4).Range.Text = "52-6="
$t.Cell(9, 5).Range.Text = "79-29="
$t.Cell(10, 1).Range.Text = "75+10="
$t.Cell(10, 2).Range.Text = "68-37="
$t.Cell(10, 3).Range.Text = "42+13="
$t.Cell(10, 4).Range.Text = "62-23="
$t.Cell(10, 5).Range.Text = "40-24="
$t.Cell(11, 1).Range.Text = "17+17="
$t.Cell(11, 2).Range.Text = "5+21="
$t.Cell(11, 3).Range.Text = "54+0="
$t.Cell(11, 4).Range.Text = "22-4="
$t.Cell(11, 5).Range.Text = "91-7="
$t.Cell(12, 1).Range.Text = "26-3="
$t.Cell(12, 2).Range.Text = "61-6="
$t.Cell(12, 3).Range.Text = "28+22="
$t.Cell(12, 4).Range.Text = "68-62="
$t.Cell(12, 5).Range.Text = "54-2="
$t.Cell(13, 1).Range.Text = "56+36="
$t.Cell(13, 2).Range.Text = "82-36="
$t.Cell(13, 3).Range.Text = "35-3="
$t.Cell(13, 4).Range.Text = "13+33="
$t.Cell(13, 5).Range.Text = "83-13="
$t.Cell(14, 1).Range.Text = "62+26="
$t.Cell(14, 2).Range.Text = "37+51="
$t.Cell(14, 3).Range.Text = "99-63="
$t.Cell(14, 4).Range.Text = "26-8="
$t.Cell(14, 5).Range.Text = "42+44="
$t.Cell(15, 1).Range.Text = "97-36="
$t.Cell(15, 2).Range.Text = "46+48="
$t.Cell(15, 3).Range.Text = "96-93="
$t.Cell(15, 4).Range.Text = "50+7="
$t.Cell(15, 5).Range.Text = "10+76="
$t.Cell(16, 1).Range.Text = "36-29="
$t.Cell(16, 2).Range.Text = "78-30="
$t.Cell(16, 3).Range.Text = "24+61="
$t.Cell(16, 4).Range.Text = "20+49="
$t.Cell(16, 5).Range.Text = "16+24="
$t.Cell(17, 1).Range.Text = "65+27="
$t.Cell(17, 2).Range.Text = "1+4="
$t.Cell(17, 3).Range.Text = "24+59="
$t.Cell(17, 4).Range.Text = "49+43="
$t.Cell(17, 5).Range.Text = "94-35="
$t.Cell(18, 1).Range.Text = "0+96="
$t.Cell(18, 2).Range.Text = "63+0="
$t.Cell(18, 3).Range.Text = "35-19="
$t.Cell(18, 4).Range.Text = "83-74="
$t.Cell(18, 5).Range.Text = "84-17="
$t.Cell(19, 1).Range.Text = "10+67="
$t.Cell(19, 2).Range.Text = "81-47="
$t.Cell(19, 3).Range.Text = "63+22="
$t.Cell(19, 4).Range.Text = "42+53="
$t.Cell(19, 5).Range.Text = "17+65="
$t.Cell(20, 1).Range.Text = "4+60="
$t.Cell(20, 2).Range.Text = "39-11="
$t.Cell(20, 3).Range.Text = "88-12="
$t.Cell(20, 4).Range.Text = "82-0="
$t.Cell(20, 5).Range.Text = "28+7="
